$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Rows(14).Insert()
$ws.Range("A14").Value = "thawing_timer"
$ws.Range("B14").Value = "Thawing Time"
[void]$ws.Range("B14").Select()
